# Apply the cryptos-list refresh described in the commit:
#   "Updated cryptos list on Sun Jul 21 21:20:48 UTC 2024 with GitHub Actions"
#
# Coin name / link cells (columns B & C) and price/volume cells (columns D & E)
# are plain text in the workbook (coinranking.com prices are stored as text,
# e.g. thousand separators are literal dots like "67.741.12"). Most new values
# can be written directly. A handful of column-D prices are plain decimal
# numbers (e.g. "599.28"); for those we prefix with a leading apostrophe so
# Excel keeps storing them as text instead of silently converting them to the
# Number type, matching the source data format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    D2 = '67.741.12'
    E2 = '  +1.04%  '
    D3 = '3.497.19'
    E3 = '  -0.40%  '
    E4 = '  -0.02%  '
    E5 = '  +0.73%  '
    E6 = '  +4.78%  '
    E7 = '  +0.02%  '
    E8 = '  +0.94%  '
    E9 = '  +4.64%  '
    E10 = '  -2.37%  '
    E11 = '  -0.09%  '
    D12 = '4.103.13'
    E12 = '  -0.36%  '
    E13 = '  +12.83%  '
    E14 = '  -0.18%  '
    D15 = '67.686.78'
    E15 = '  +0.99%  '
    E16 = '  +0.62%  '
    D17 = '3.488.95'
    E17 = '  -0.36%  '
    E19 = '  +3.09%  '
    E20 = '  +0.18%  '
    E21 = '  +1.31%  '
    E22 = '  +0.28%  '
    E23 = '  +1.04%  '
    E24 = '  +0.16%  '
    B25 = 'LEO'
    C25 = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    E25 = '  +0.23%  '
    B26 = 'PEPE'
    C26 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
    E26 = '  +2.30%  '
    E27 = '  +1.92%  '
    E28 = '  -0.78%  '
    E29 = '  -0.20%  '
    E30 = '  +0.17%  '
    E31 = '  +0.66%  '
    E32 = '  +0.22%  '
    E33 = '  +0.25%  '
    E34 = '  +0.47%  '
    E35 = '  +0.10%  '
    E37 = '  +0.54%  '
    E38 = '  +2.57%  '
    E39 = '  -2.37%  '
    E40 = '  +2.89%  '
    E41 = '  +1.14%  '
    E42 = '  +1.77%  '
    E43 = '  +2.44%  '
    E44 = '  +1.05%  '
    E45 = '  -0.93%  '
    D46 = '2.821.45'
    E46 = '  +0.96%  '
    E47 = '  -1.22%  '
    E48 = '  -0.24%  '
    E49 = '  +1.31%  '
    E50 = '  -1.17%  '
    E51 = '  +1.00%  '
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$textProtectedUpdates = @{
    D5 = '599.28'
    D6 = '182.29'
    D10 = '7.10'
    D11 = '0.434'
    D13 = '32.46'
    D18 = '6.38'
    D19 = '14.71'
    D20 = '395.96'
    D21 = '8.07'
    D22 = '73.40'
    D23 = '0.545'
    D25 = '5.70'
    D26 = '0.0000125'
    D27 = '10.44'
    D28 = '0.179'
    D29 = '0.996'
    D30 = '6.28'
    D31 = '1.46'
    D33 = '24.02'
    D37 = '164.29'
    D39 = '0.873'
    D41 = '4.75'
    D42 = '27.73'
    D44 = '26.63'
    D47 = '42.32'
    D48 = '0.0305'
    D49 = '345.46'
    D50 = '1.08'
    D51 = '33.60'
}
foreach ($ref in $textProtectedUpdates.Keys) {
    $ws.Range($ref).Value = "'" + $textProtectedUpdates[$ref]
}
